# This workbook tracks weekly "Zapallo italiano" price records for the
# "Terminal Hortofrutícola Agro Chillán" market. The edit adds one new
# weekly record, inserted as a new row 72 (pushing every subsequent
# record down by one row, from 72-117 to 73-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 72; existing rows 72..117 shift down
# to 73..118 (dimension grows from A1:R117 to A1:R118).
$ws.Range("A72").EntireRow.Insert()

# Populate the newly inserted row 72 with the new weekly record.
$ws.Cells.Item(72, 1).Value2  = 7
$ws.Cells.Item(72, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value2  = "Ñuble"
$ws.Cells.Item(72, 4).Value2  = 44438
$ws.Cells.Item(72, 5).Value2  = 16
$ws.Cells.Item(72, 6).Value2  = 100112032
$ws.Cells.Item(72, 7).Value2  = "Zapallo italiano"
$ws.Cells.Item(72, 8).Value2  = "Sin especificar"
$ws.Cells.Item(72, 9).Value2  = "Primera"
$ws.Cells.Item(72, 10).Value2 = 160
$ws.Cells.Item(72, 11).Value2 = 14000
$ws.Cells.Item(72, 12).Value2 = 15000
$ws.Cells.Item(72, 13).Value2 = 14500
$ws.Cells.Item(72, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item(72, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value2 = 290
$ws.Cells.Item(72, 17).Value2 = 50
$ws.Cells.Item(72, 18).Value2 = "Hortaliza"
